# Update the date line and the multiplication problems in the table
# as described by the commit diff.

$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $true, $false, $false, $false, `
                             $true, 1, $false, $new, 2)
}

Replace-Text "2024-04-07 Sunday" "2024-04-08 Monday"

Replace-Text "316×4=" "487×2="
Replace-Text "190×7=" "493×9="
Replace-Text "152×5=" "918×4="
Replace-Text "129×6=" "496×5="
Replace-Text "212×2=" "368×3="
Replace-Text "959×8=" "360×8="
Replace-Text "514×9=" "903×8="
Replace-Text "261×4=" "744×9="
Replace-Text "521×4=" "972×6="
Replace-Text "120×3=" "886×8="
Replace-Text "243×2=" "793×8="
Replace-Text "711×6=" "623×2="
Replace-Text "412×4=" "633×6="
Replace-Text "670×4=" "969×5="
Replace-Text "200×7=" "380×2="
Replace-Text "967×3=" "439×2="
Replace-Text "921×3=" "787×2="
Replace-Text "900×5=" "900×6="
Replace-Text "550×5=" "751×5="
Replace-Text "174×6=" "432×4="
Replace-Text "332×7=" "992×9="
Replace-Text "673×4=" "361×4="
Replace-Text "403×4=" "807×3="
Replace-Text "991×6=" "484×8="
Replace-Text "663×5=" "283×8="

$d.Save()
